$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'279.12"
$ws.Range("E2").Value = "'6.48%"
$ws.Range("D3").Value = "'27.10"
$ws.Range("E3").Value = "'1.52%"
$ws.Range("D4").Value = "'4.901"
$ws.Range("E4").Value = "'4.43%"
$ws.Range("E5").Value = "'3.87%"
$ws.Range("D6").Value = "'6.945"
$ws.Range("E6").Value = "'3.64%"
$ws.Range("D7").Value = "'3.354"
$ws.Range("E7").Value = "'5.97%"
$ws.Range("D8").Value = "'0.8821"
$ws.Range("E8").Value = "'3.75%"
$ws.Range("E9").Value = "'4.15%"
$ws.Range("D10").Value = "'0.1470"
$ws.Range("E10").Value = "'4.68%"
$ws.Range("D11").Value = "'0.05154"
$ws.Range("E11").Value = "'2.12%"
$ws.Range("D12").Value = "'0.07418"
$ws.Range("E12").Value = "'4.36%"
$ws.Range("D13").Value = "'0.03146"
$ws.Range("E13").Value = "'-0.07%"
$ws.Range("D14").Value = "'0.09073"
$ws.Range("E14").Value = "'0.26%"
$ws.Range("D15").Value = "'0.001559"
$ws.Range("E15").Value = "'1.67%"
$ws.Range("D16").Value = "'0.0006265"
$ws.Range("E16").Value = "'1.27%"
$ws.Range("D17").Value = "'0.005852"
$ws.Range("E17").Value = "'-2.88%"
$ws.Range("D18").Value = "'3.478"
$ws.Range("E18").Value = "'0.75%"
$ws.Range("D19").Value = "'2.296"
$ws.Range("E19").Value = "'5.97%"
$ws.Range("D21").Value = "'0.1340"
$ws.Range("E21").Value = "'3.03%"
$ws.Range("D22").Value = "'3.893"
$ws.Range("E22").Value = "'-4.67%"
$ws.Range("D23").Value = "'0.04322"
$ws.Range("E23").Value = "'1.94%"
$ws.Range("D24").Value = "'0.001175"
$ws.Range("E24").Value = "'-0.11%"
$ws.Range("D25").Value = "'0.003611"
$ws.Range("E25").Value = "'-11.01%"
$ws.Range("E26").Value = "'-0.12%"
$ws.Range("D27").Value = "'0.0001695"
$ws.Range("E27").Value = "'-12.52%"
$ws.Range("D40").Value = "'0.04044"
$ws.Range("E40").Value = "'2.61%"
$ws.Range("D41").Value = "'0.006615"
$ws.Range("E41").Value = "'58.26%"
$ws.Range("D42").Value = "'0.1163"
$ws.Range("E42").Value = "'4.67%"
$ws.Range("D43").Value = "'0.002338"
$ws.Range("E43").Value = "'10.24%"
$ws.Range("D44").Value = "'0.01250"
$ws.Range("E44").Value = "'8.81%"
$ws.Range("D45").Value = "'0.00005225"
$ws.Range("E45").Value = "'2.41%"
$ws.Range("E46").Value = "'0.00%"
$ws.Range("E47").Value = "'823.55%"
$ws.Range("D48").Value = "'0.02251"
$ws.Range("E48").Value = "'6.12%"
$ws.Range("E49").Value = "'0.00%"
$ws.Range("E50").Value = "'-0.07%"
